$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primjer 4.1")

# Weighted-sum ("Bodovi") column for each expert/row (A, B, C)
$ws.Range("G6").Formula = "=C6*C5+D6*D5+E6*E5"
$ws.Range("G7").Formula = "=C7*C5+D7*D5+E7*E5"
$ws.Range("G8").Formula = "=C8*C5+D8*D5+E8*E5"

# New row 11: "Bodovi" label + weighting values
$ws.Range("B11").Value = "Bodovi"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0

$ws.Range("B4").Copy()
$ws.Range("B11:E11").PasteSpecial(-4122)
$ws.Range("B11:E11").HorizontalAlignment = -4108
$ws.Range("C11:E11").Interior.Pattern = -4142

$ws.Range("G9").Select()
